# G8.C3.W2 PPTX slide layout fixes
# Slide 3: reposition the "Driving Question" text box, the white separator
#          line beneath it, the "And how does this connect..." text, and
#          the mission call-out box/text so nothing overlaps.
# Slide 4: shrink the title font from 24pt to 20pt so it fits on one line.

$p = $ppt.ActivePresentation

# ---- Slide 3 -------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "Why do whale flippers..." / "What does this tell us about whale
# ancestors?" question text box - move up slightly and shrink its height
# so the separator line below sits closer to the text.
# (the tiny +0.00001 nudges avoid single-precision round-trip truncation
# in the host's points->EMU conversion landing one EMU short of target)
$questionBox = $s3.Shapes.Item(3)
$questionBox.Top = 93.60001
$questionBox.Height = 79.20001

# White rounded-rectangle divider line - move up to sit right below the
# question text instead of floating further down the slide.
$dividerLine = $s3.Shapes.Item(4)
$dividerLine.Top = 187.20001

# "And how does this connect to natural selection from Week 1?" text box
$connectText = $s3.Shapes.Item(5)
$connectText.Top = 201.6

# Teal "Your Mission" rounded-rectangle background
$missionBox = $s3.Shapes.Item(6)
$missionBox.Top = 259.2

# "Your Mission: Use anatomical & fossil evidence..." text
$missionText = $s3.Shapes.Item(7)
$missionText.Top = 266.40001

# ---- Slide 4 -------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Title "What You Already Know (Week 1 Natural Selection)" - reduce font
# size from 24pt to 20pt so it fits on a single line and no longer
# overlaps the "Variation exists" / "Selection pressure" boxes below it.
$title = $s4.Shapes.Item(1)
$title.TextFrame.TextRange.Font.Size = 20
